# Update the "build timestamp" embedded in the version strings across the workbook.
# Old build timestamp: January 30 2026 16.19.47 EST
# New build timestamp: February 02 2026 12.49.33 EST

$wb = $excel.ActiveWorkbook

$newStamp = "February 02 2026 12.49.33 EST"

# --- Sheet "About" ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: mines - January 30 (built on " + $newStamp + ")"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Tikhova Coal Mine, Russia, M0856, version 'mines - January 30 (built on " + $newStamp + ")'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet "Boundaries and methane sources" ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 12; $row++) {
    $cell = $wsData.Range("S" + $row)
    $cell.Value = "mines - January 30 (built on " + $newStamp + ")"
}
